# CCC19 Derived Variables Spreadsheet - add three new derived variables
# (X4/quality, X5/problems, X6/ccc19cci) to the Table1 listobject on Sheet1.
# This mirrors the fix described in the commit message: the new rows feed
# the quality-score / composite variables that relied on the corrected
# O2-requirement baseline variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$row1 = $lo.ListRows.Add()
$row1.Range.Item(1,1).Value = "X4"
$row1.Range.Item(1,2).Value = "quality"
$row1.Range.Item(1,3).Value = "Other"
$row1.Range.Item(1,4).Value = "Quality score"

$row2 = $lo.ListRows.Add()
$row2.Range.Item(1,1).Value = "X5"
$row2.Range.Item(1,2).Value = "problems"
$row2.Range.Item(1,3).Value = "Other"
$row2.Range.Item(1,4).Value = "Enumerated problems that feed quality score"

$row3 = $lo.ListRows.Add()
$row3.Range.Item(1,1).Value = "X6"
$row3.Range.Item(1,2).Value = "ccc19cci"
$row3.Range.Item(1,3).Value = "Other"
$row3.Range.Item(1,4).Value = "CCC19 modified Charlson"

# Reflect the author's final cursor/selection position (cell D107) like in
# the saved workbook's sheetView.
$ws.Range("D107").Select()
